$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 2.430374
$ws.Range("H2").Value2 = 7.291122000000001
$ws.Range("I2").Value2 = 0.009222757332915244
$ws.Range("J2").Value2 = 0.009222757332915246
$ws.Range("M2").Value2 = 91.74689966666665
$ws.Range("N2").Value2 = 275.2406989999999
$ws.Range("O2").Value2 = 0.1908387282982634
$ws.Range("P2").Value2 = 0.1908387282982634
$ws.Range("Q2").Value2 = 222.9792795304753
$ws.Range("R2").Value2 = 2006.813515774278
$ws.Range("S2").Value2 = 0.001760059280817028
$ws.Range("T2").Value2 = 0.001760059280817029

$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 2.430374
$ws.Range("H3").Value2 = 7.291122000000001
$ws.Range("I3").Value2 = 0.009222757332915244
$ws.Range("J3").Value2 = 0.009222757332915246
$ws.Range("O3").Value2 = 0.296899627499751
$ws.Range("P3").Value2 = 0.296899627499751
$ws.Range("Q3").Value2 = 346.9026733886667
$ws.Range("R3").Value2 = 3122.124060498
$ws.Range("S3").Value2 = 0.002738233216663132
$ws.Range("T3").Value2 = 0.002738233216663133

$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 2.430374
$ws.Range("H4").Value2 = 7.291122000000001
$ws.Range("I4").Value2 = 0.009222757332915244
$ws.Range("J4").Value2 = 0.009222757332915246
$ws.Range("M4").Value2 = 167.6324513333334
$ws.Range("N4").Value2 = 502.8973540000001
$ws.Range("O4").Value2 = 0.348684957750095
$ws.Range("P4").Value2 = 0.348684957750095
$ws.Range("Q4").Value2 = 407.4095512767988
$ws.Range("R4").Value2 = 3666.685961491189
$ws.Range("S4").Value2 = 0.003215836750966931
$ws.Range("T4").Value2 = 0.003215836750966931

$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 2.430374
$ws.Range("H5").Value2 = 7.291122000000001
$ws.Range("I5").Value2 = 0.009222757332915244
$ws.Range("J5").Value2 = 0.009222757332915246
$ws.Range("M5").Value2 = 78.64050433333334
$ws.Range("N5").Value2 = 235.921513
$ws.Range("O5").Value2 = 0.1635766864518907
$ws.Range("P5").Value2 = 0.1635766864518907
$ws.Range("Q5").Value2 = 191.1258370786207
$ws.Range("R5").Value2 = 1720.132533707586
$ws.Range("S5").Value2 = 0.001508628084468152
$ws.Range("T5").Value2 = 0.001508628084468153

$ws.Range("I6").Value2 = 0.5480399755605952
$ws.Range("J6").Value2 = 0.5480399755605954
$ws.Range("M6").Value2 = 91.74689966666665
$ws.Range("N6").Value2 = 275.2406989999999
$ws.Range("O6").Value2 = 0.1908387282982634
$ws.Range("P6").Value2 = 0.1908387282982634
$ws.Range("Q6").Value2 = 13250.0026286362
$ws.Range("R6").Value2 = 119250.0236577258
$ws.Range("S6").Value2 = 0.1045872519925953
$ws.Range("T6").Value2 = 0.1045872519925954

$ws.Range("I7").Value2 = 0.5480399755605952
$ws.Range("J7").Value2 = 0.5480399755605954
$ws.Range("O7").Value2 = 0.296899627499751
$ws.Range("P7").Value2 = 0.296899627499751
$ws.Range("S7").Value2 = 0.1627128645989133
$ws.Range("T7").Value2 = 0.1627128645989134

$ws.Range("I8").Value2 = 0.5480399755605952
$ws.Range("J8").Value2 = 0.5480399755605954
$ws.Range("M8").Value2 = 167.6324513333334
$ws.Range("N8").Value2 = 502.8973540000001
$ws.Range("O8").Value2 = 0.348684957750095
$ws.Range("P8").Value2 = 0.348684957750095
$ws.Range("Q8").Value2 = 24209.32400856239
$ws.Range("R8").Value2 = 217883.9160770615
$ws.Range("S8").Value2 = 0.1910932957237093
$ws.Range("T8").Value2 = 0.1910932957237093

$ws.Range("I9").Value2 = 0.5480399755605952
$ws.Range("J9").Value2 = 0.5480399755605954
$ws.Range("M9").Value2 = 78.64050433333334
$ws.Range("N9").Value2 = 235.921513
$ws.Range("O9").Value2 = 0.1635766864518907
$ws.Range("P9").Value2 = 0.1635766864518907
$ws.Range("Q9").Value2 = 11357.1891030615
$ws.Range("R9").Value2 = 102214.7019275535
$ws.Range("S9").Value2 = 0.08964656324537733
$ws.Range("T9").Value2 = 0.08964656324537734

$ws.Range("G10").Value2 = 116.470388
$ws.Range("H10").Value2 = 349.411164
$ws.Range("I10").Value2 = 0.4419805861132828
$ws.Range("J10").Value2 = 0.4419805861132828
$ws.Range("M10").Value2 = 91.74689966666665
$ws.Range("N10").Value2 = 275.2406989999999
$ws.Range("O10").Value2 = 0.1908387282982634
$ws.Range("P10").Value2 = 0.1908387282982634
$ws.Range("Q10").Value2 = 10685.79700197373
$ws.Range("R10").Value2 = 96172.17301776362
$ws.Range("S10").Value2 = 0.08434701298637998
$ws.Range("T10").Value2 = 0.08434701298637999

$ws.Range("G11").Value2 = 116.470388
$ws.Range("H11").Value2 = 349.411164
$ws.Range("I11").Value2 = 0.4419805861132828
$ws.Range("J11").Value2 = 0.4419805861132828
$ws.Range("O11").Value2 = 0.296899627499751
$ws.Range("P11").Value2 = 0.296899627499751
$ws.Range("Q11").Value2 = 16624.55612503067
$ws.Range("R11").Value2 = 149621.005125276
$ws.Range("S11").Value2 = 0.1312238713791553
$ws.Range("T11").Value2 = 0.1312238713791553

$ws.Range("G12").Value2 = 116.470388
$ws.Range("H12").Value2 = 349.411164
$ws.Range("I12").Value2 = 0.4419805861132828
$ws.Range("J12").Value2 = 0.4419805861132828
$ws.Range("M12").Value2 = 167.6324513333334
$ws.Range("N12").Value2 = 502.8973540000001
$ws.Range("O12").Value2 = 0.348684957750095
$ws.Range("P12").Value2 = 0.348684957750095
$ws.Range("Q12").Value2 = 19524.21664818445
$ws.Range("R12").Value2 = 175717.9498336601
$ws.Range("S12").Value2 = 0.1541119819952723
$ws.Range("T12").Value2 = 0.1541119819952723

$ws.Range("G13").Value2 = 116.470388
$ws.Range("H13").Value2 = 349.411164
$ws.Range("I13").Value2 = 0.4419805861132828
$ws.Range("J13").Value2 = 0.4419805861132828
$ws.Range("M13").Value2 = 78.64050433333334
$ws.Range("N13").Value2 = 235.921513
$ws.Range("O13").Value2 = 0.1635766864518907
$ws.Range("P13").Value2 = 0.1635766864518907
$ws.Range("Q13").Value2 = 9159.290052219016
$ws.Range("R13").Value2 = 82433.61046997113
$ws.Range("S13").Value2 = 0.07229771975247534
$ws.Range("T13").Value2 = 0.07229771975247534

$ws.Range("E14").Value2 = 3
$ws.Range("F14").Value2 = 1
$ws.Range("G14").Value2 = 0.1994
$ws.Range("H14").Value2 = 0.5982000000000001
$ws.Range("I14").Value2 = 0.0007566809932065188
$ws.Range("J14").Value2 = 0.0007566809932065189
$ws.Range("M14").Value2 = 91.74689966666665
$ws.Range("N14").Value2 = 275.2406989999999
$ws.Range("O14").Value2 = 0.1908387282982634
$ws.Range("P14").Value2 = 0.1908387282982634
$ws.Range("Q14").Value2 = 18.29433179353333
$ws.Range("R14").Value2 = 164.6489861418
$ws.Range("S14").Value2 = 0.0001444040384709989
$ws.Range("T14").Value2 = 0.000144404038470999

$ws.Range("E15").Value2 = 3
$ws.Range("F15").Value2 = 1
$ws.Range("G15").Value2 = 0.1994
$ws.Range("H15").Value2 = 0.5982000000000001
$ws.Range("I15").Value2 = 0.0007566809932065188
$ws.Range("J15").Value2 = 0.0007566809932065189
$ws.Range("O15").Value2 = 0.296899627499751
$ws.Range("P15").Value2 = 0.296899627499751
$ws.Range("Q15").Value2 = 28.46162486666667
$ws.Range("R15").Value2 = 256.1546238
$ws.Range("S15").Value2 = 0.000224658305019157
$ws.Range("T15").Value2 = 0.000224658305019157

$ws.Range("E16").Value2 = 3
$ws.Range("F16").Value2 = 1
$ws.Range("G16").Value2 = 0.1994
$ws.Range("H16").Value2 = 0.5982000000000001
$ws.Range("I16").Value2 = 0.0007566809932065188
$ws.Range("J16").Value2 = 0.0007566809932065189
$ws.Range("M16").Value2 = 167.6324513333334
$ws.Range("N16").Value2 = 502.8973540000001
$ws.Range("O16").Value2 = 0.348684957750095
$ws.Range("P16").Value2 = 0.348684957750095
$ws.Range("Q16").Value2 = 33.42591079586668
$ws.Range("R16").Value2 = 300.8331971628
$ws.Range("S16").Value2 = 0.0002638432801465149
$ws.Range("T16").Value2 = 0.000263843280146515

$ws.Range("E17").Value2 = 3
$ws.Range("F17").Value2 = 1
$ws.Range("G17").Value2 = 0.1994
$ws.Range("H17").Value2 = 0.5982000000000001
$ws.Range("I17").Value2 = 0.0007566809932065188
$ws.Range("J17").Value2 = 0.0007566809932065189
$ws.Range("M17").Value2 = 78.64050433333334
$ws.Range("N17").Value2 = 235.921513
$ws.Range("O17").Value2 = 0.1635766864518907
$ws.Range("P17").Value2 = 0.1635766864518907
$ws.Range("Q17").Value2 = 15.68091656406667
$ws.Range("R17").Value2 = 141.1282490766
$ws.Range("S17").Value2 = 0.000123775369569848
$ws.Range("T17").Value2 = 0.000123775369569848
